# Auto-generated: applies numeric cell updates across all 8 sheets
# to match the target diff (market-price-derived profit columns H-N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 470
$ws.Range("I2").Value = 300
$ws.Range("K2").Value = 300
$ws.Range("M2").Value = -187
$ws.Range("H17").Value = 1070.1466
$ws.Range("J17").Value = 1070.1466
$ws.Range("L17").Value = 3210.4398
$ws.Range("N17").Value = -3546.4398
$ws.Range("H132").Value = 1161.881
$ws.Range("I132").Value = 984.0526
$ws.Range("K132").Value = 2952.1578
$ws.Range("M132").Value = -422.1578
$ws.Range("H137").Value = 1358.5
$ws.Range("I137").Value = 1202.2632
$ws.Range("K137").Value = 3606.7896
$ws.Range("M137").Value = -1056.7896
$ws.Range("H138").Value = 3640.889
$ws.Range("I138").Value = 3478.5293
$ws.Range("J138").Value = 3916.9
$ws.Range("K138").Value = 10435.5879
$ws.Range("L138").Value = 11750.7
$ws.Range("M138").Value = -5295.5879
$ws.Range("N138").Value = -22030.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2623.716
$ws.Range("I32").Value = 1885.0588
$ws.Range("K32").Value = 1885.0588
$ws.Range("M32").Value = -1598.0588
$ws.Range("H61").Value = 2069.0667
$ws.Range("I61").Value = 1142.9546
$ws.Range("K61").Value = 1142.9546
$ws.Range("M61").Value = -930.9546
$ws.Range("H97").Value = 831
$ws.Range("I97").Value = 608
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 608
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -112
$ws.Range("N97").Value = -2492
$ws.Range("H109").Value = 46451.332
$ws.Range("J109").Value = 46451.332
$ws.Range("L109").Value = 46451.332
$ws.Range("N109").Value = -49225.332
$ws.Range("H136").Value = 2069.0667
$ws.Range("I136").Value = 1142.9546
$ws.Range("K136").Value = 3428.8638
$ws.Range("M136").Value = -878.8638000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 57029
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("H99").Value = 1420
$ws.Range("I99").Value = 1209.091
$ws.Range("K99").Value = 1209.091
$ws.Range("M99").Value = 288.9090000000001
$ws.Range("H108").Value = 94985.5
$ws.Range("J108").Value = 94985.5
$ws.Range("L108").Value = 94985.5
$ws.Range("N108").Value = -102665.5
$ws.Range("M37").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1061996.9
$ws.Range("I58").Value = 1553923.9
$ws.Range("J58").Value = 2461.6155
$ws.Range("K58").Value = 1553923.9
$ws.Range("L58").Value = 2461.6155
$ws.Range("M58").Value = -1553720.9
$ws.Range("N58").Value = -2867.6155
$ws.Range("H132").Value = 1845.6531
$ws.Range("I132").Value = 1245.3226
$ws.Range("J132").Value = 2879.5557
$ws.Range("K132").Value = 3735.9678
$ws.Range("L132").Value = 8638.667099999999
$ws.Range("M132").Value = -1205.9678
$ws.Range("N132").Value = -13698.6671
$ws.Range("H136").Value = 1061996.9
$ws.Range("I136").Value = 1553923.9
$ws.Range("J136").Value = 2461.6155
$ws.Range("K136").Value = 4661771.699999999
$ws.Range("L136").Value = 7384.8465
$ws.Range("M136").Value = -4659221.699999999
$ws.Range("N136").Value = -12484.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1332.6666
$ws.Range("J80").Value = 1332.6666
$ws.Range("L80").Value = 3997.9998
$ws.Range("N80").Value = -5869.9998
$ws.Range("H83").Value = 1332.6666
$ws.Range("J83").Value = 1332.6666
$ws.Range("L83").Value = 11993.9994
$ws.Range("N83").Value = -21353.9994
$ws.Range("H105").Value = 2816.7646
$ws.Range("J105").Value = 2937.3125
$ws.Range("L105").Value = 8811.9375
$ws.Range("N105").Value = -14053.9375
$ws.Range("H107").Value = 679.9167
$ws.Range("I107").Value = 281
$ws.Range("J107").Value = 812.8889
$ws.Range("K107").Value = 843
$ws.Range("L107").Value = 2438.6667
$ws.Range("M107").Value = 1077
$ws.Range("N107").Value = -6278.6667
$ws.Range("H122").Value = 1009.1818
$ws.Range("J122").Value = 1296.8334
$ws.Range("L122").Value = 11671.5006
$ws.Range("N122").Value = -16571.5006
$ws.Range("H131").Value = 823.45
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 823.45
$ws.Range("K131").Value = 0
$ws.Range("N131").Value = -12550.35
$ws.Range("M131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 62.285713
$ws.Range("I2").Value = 19.666666
$ws.Range("J2").Value = 73.90909000000001
$ws.Range("K2").Value = 19.666666
$ws.Range("L2").Value = 73.90909000000001
$ws.Range("M2").Value = 93.33333400000001
$ws.Range("N2").Value = -299.90909
$ws.Range("H15").Value = 14499.5
$ws.Range("J15").Value = 14499.5
$ws.Range("L15").Value = 14499.5
$ws.Range("N15").Value = -15075.5
$ws.Range("H40").Value = 57018
$ws.Range("J40").Value = 57018
$ws.Range("L40").Value = 57018
$ws.Range("N40").Value = -57320
$ws.Range("H81").Value = 14499.5
$ws.Range("J81").Value = 14499.5
$ws.Range("L81").Value = 14499.5
$ws.Range("N81").Value = -16495.5
$ws.Range("H84").Value = 14499.5
$ws.Range("J84").Value = 14499.5
$ws.Range("L84").Value = 43498.5
$ws.Range("N84").Value = -53482.5
$ws.Range("H122").Value = 1513.7931
$ws.Range("I122").Value = 1222.9286
$ws.Range("J122").Value = 1785.2667
$ws.Range("K122").Value = 3668.7858
$ws.Range("L122").Value = 5355.800099999999
$ws.Range("M122").Value = -1218.7858
$ws.Range("N122").Value = -10255.8001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2259.3076
$ws.Range("I46").Value = 1283
$ws.Range("J46").Value = 3096.1428
$ws.Range("K46").Value = 1283
$ws.Range("L46").Value = 3096.1428
$ws.Range("M46").Value = -1095
$ws.Range("N46").Value = -3472.1428
$ws.Range("H94").Value = 54164.5
$ws.Range("J94").Value = 54164.5
$ws.Range("L94").Value = 54164.5
$ws.Range("N94").Value = -55516.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 57062
$ws.Range("J38").Value = 57062
$ws.Range("L38").Value = 57062
$ws.Range("N38").Value = -58008
$ws.Range("H100").Value = 541.8570999999999
$ws.Range("I100").Value = 407.81818
$ws.Range("J100").Value = 1033.3334
$ws.Range("K100").Value = 815.63636
$ws.Range("L100").Value = 2066.6668
$ws.Range("M100").Value = -274.63636
$ws.Range("N100").Value = -3148.6668
$ws.Range("H122").Value = 88311.11
$ws.Range("I122").Value = 195644
$ws.Range("J122").Value = 2444.8
$ws.Range("K122").Value = 586932
$ws.Range("L122").Value = 7334.400000000001
$ws.Range("M122").Value = -584482
$ws.Range("N122").Value = -12234.4
$ws.Range("H132").Value = 1099.814
$ws.Range("I132").Value = 796.5625
$ws.Range("J132").Value = 1982
$ws.Range("K132").Value = 2389.6875
$ws.Range("L132").Value = 5946
$ws.Range("M132").Value = 140.3125
$ws.Range("N132").Value = -11006

Write-Host "Applied all cell updates."